$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00009552326474482342
$ws.Range("C2").Value = 1.62698769954209399
$ws.Range("D2").Value = 18.71679738969934093
$ws.Range("E2").Value = 2797.56581773474408692
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2817.90969834725001419
